# #5: fund, bonds, otherbonds, antique done
#
# Sheet "具有相當價值之財產" (property of considerable value) gains a
# standard trailer of metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that already exist on
# the other sheets (存款/股票/...). A brand-new shared string "otherbonds" is
# used as the property_category for every data row on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

# ---- Header row (row 1): columns B..E get re-labelled, F..L are new ----
# Copy formatting (bold + border, same xf as the existing header cells)
# from the existing E1 header cell onto the newly-appended header cells.
$ws.Range("E1").Copy()
$ws.Range("F1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "quantity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "property_category"
$ws.Cells.Item(1, 7).Value = "category"
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"
$ws.Cells.Item(1, 11).Value = "source_file"
$ws.Cells.Item(1, 12).Value = "index"

# ---- Data rows 2..19: append F..L, same value on every row except the ----
# ---- trailing "index" column which mirrors column A on that row.      ----
$firstDataRow = 2
$lastDataRow = 19

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    # Copy formatting from the row's own column E cell (same border/font as
    # the rest of the row) onto the new trailing cells.
    $ws.Range("E$r").Copy()
    $ws.Range("F$r`:L$r").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $indexValue = $ws.Cells.Item($r, 1).Value2

    $ws.Cells.Item($r, 6).Value = "otherbonds"
    $ws.Cells.Item($r, 7).Value = "normal"
    # Quote-prefix so the engine keeps this as literal text instead of
    # re-parsing "2012-04-19" into a date serial number.
    $ws.Cells.Item($r, 8).Value = "'2012-04-19"
    $ws.Cells.Item($r, 9).Value = "張慶忠"
    $ws.Cells.Item($r, 10).Value = 1347
    $ws.Cells.Item($r, 11).Value = "tmp93201"
    $ws.Cells.Item($r, 12).Value = $indexValue
}
